$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '75.622.57'
$ws.Range("E2").Value = '  +9.02%  '

$ws.Range("D3").Value = '2.705.60'
$ws.Range("E3").Value = '  +11.53%  '

$ws.Range("E4").Value = '  -0.08%  '

Set-TextValue 'D5' '189.71'
$ws.Range("E5").Value = '  +14.12%  '

Set-TextValue 'D6' '590.30'
$ws.Range("E6").Value = '  +4.84%  '

$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +5.64%  '

Set-TextValue 'D9' '0.198'
$ws.Range("E9").Value = '  +16.92%  '

$ws.Range("D10").Value = '2.701.04'
$ws.Range("E10").Value = '  +11.37%  '

Set-TextValue 'D11' '0.162'
$ws.Range("E11").Value = '  +1.36%  '

Set-TextValue 'D12' '0.361'
$ws.Range("E12").Value = '  +7.96%  '

Set-TextValue 'D13' '4.76'
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '75.626.88'
$ws.Range("E14").Value = '  +9.19%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.183.22'
$ws.Range("E15").Value = '  +10.78%  '

$ws.Range("E16").Value = '  +7.20%  '

Set-TextValue 'D17' '26.78'
$ws.Range("E17").Value = '  +11.95%  '

$ws.Range("D18").Value = '2.700.43'
$ws.Range("E18").Value = '  +11.77%  '

Set-TextValue 'D19' '9.50'
$ws.Range("E19").Value = '  +33.12%  '

Set-TextValue 'D20' '12.17'
$ws.Range("E20").Value = '  +12.72%  '

Set-TextValue 'D21' '380.02'
$ws.Range("E21").Value = '  +10.77%  '

$ws.Range("E22").Value = '  +16.79%  '

Set-TextValue 'D23' '4.09'
$ws.Range("E23").Value = '  +5.74%  '

Set-TextValue 'D24' '6.29'
$ws.Range("E24").Value = '  +4.74%  '

Set-TextValue 'D25' '1.00'
$ws.Range("E25").Value = '  +0.03%  '

Set-TextValue 'D26' '70.75'
$ws.Range("E26").Value = '  +7.19%  '

Set-TextValue 'D27' '4.23'
$ws.Range("E27").Value = '  +10.76%  '

Set-TextValue 'D28' '9.56'
$ws.Range("E28").Value = '  +12.66%  '

$ws.Range("E29").Value = '  +10.91%  '

$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").Value = '0.0₃0971'
$ws.Range("E31").Value = '  +14.44%  '

Set-TextValue 'D32' '523.65'
$ws.Range("E32").Value = '  +15.56%  '

$ws.Range("E33").Value = '  +14.22%  '

Set-TextValue 'D34' '7.88'
$ws.Range("E34").Value = '  +6.65%  '

Set-TextValue 'D35' '1.78'
$ws.Range("E35").Value = '  +10.32%  '

$ws.Range("E36").Value = '  -0.20%  '

$ws.Range("E37").Value = '  +8.64%  '

Set-TextValue 'D38' '162.49'
$ws.Range("E38").Value = '  +2.20%  '

Set-TextValue 'D39' '19.43'

Set-TextValue 'D40' '19.38'
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '173.83'
$ws.Range("E42").Value = '  +28.10%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D43' '5.07'
$ws.Range("E43").Value = '  +15.12%  '

Set-TextValue 'D44' '1.72'
$ws.Range("E44").Value = '  +13.27%  '

Set-TextValue 'D45' '0.334'
$ws.Range("E45").Value = '  +10.16%  '

$ws.Range("E46").Value = '  +11.54%  '

Set-TextValue 'D47' '2.42'
$ws.Range("E47").Value = '  +15.95%  '

Set-TextValue 'D48' '39.13'
$ws.Range("E48").Value = '  +3.57%  '

Set-TextValue 'D49' '0.0851'
$ws.Range("E49").Value = '  +18.05%  '

Set-TextValue 'D50' '3.69'
$ws.Range("E50").Value = '  +8.87%  '

Set-TextValue 'D51' '0.546'
$ws.Range("E51").Value = '  +11.68%  '

Write-Host "Applied crypto list update"
